$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("debt_schedule")

# ---- Row 2 ----
$ws.Cells.Item(2,1).Value = "First Lien Notes 2026"
$ws.Cells.Item(2,2).Value = "secured"
$ws.Cells.Item(2,3).Value = 10.5
$ws.Cells.Item(2,4).Formula = "'" + "2026-04-24"
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).Value = 429.3
$ws.Cells.Item(2,6).Value = "No call provision"
$ws.Cells.Item(2,7).Value = "Senior Secured First Lien Notes due 2026, secured by substantially all assets"
$ws.Cells.Item(2,8).ClearContents()
$ws.Cells.Item(2,9).ClearContents()
$ws.Cells.Item(2,10).ClearContents()
$ws.Cells.Item(2,11).ClearContents()

# ---- Row 3 ----
$ws.Cells.Item(3,1).Value = "First Lien Notes 2029"
$ws.Cells.Item(3,2).Value = "secured"
$ws.Cells.Item(3,3).Value = 7.5
$ws.Cells.Item(3,4).Formula = "'" + "2029-02-15"
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).Value = 800
$ws.Cells.Item(3,6).Value = "Callable at 103.75% (2025)"
$ws.Cells.Item(3,7).Value = "Senior Secured First Lien Notes due 2029, callable premium declines 1.25% annually"
$ws.Cells.Item(3,8).ClearContents()
$ws.Cells.Item(3,9).ClearContents()
$ws.Cells.Item(3,10).ClearContents()
$ws.Cells.Item(3,11).ClearContents()

# ---- Row 4 ----
$ws.Cells.Item(4,1).Value = "Second Lien Notes 2026"
$ws.Cells.Item(4,2).Value = "unsecured"
$ws.Cells.Item(4,3).Value = 10
$ws.Cells.Item(4,4).Formula = "'" + "2026-06-15"
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,5).Value = 199.2
$ws.Cells.Item(4,6).Value = "No call provision"
$ws.Cells.Item(4,7).Value = "Senior Unsecured Second Lien Notes due 2026"
$ws.Cells.Item(4,8).ClearContents()
$ws.Cells.Item(4,9).ClearContents()
$ws.Cells.Item(4,10).ClearContents()
$ws.Cells.Item(4,11).ClearContents()

# ---- Row 5 ----
$ws.Cells.Item(5,1).Value = "Convertible Notes 2026"
$ws.Cells.Item(5,2).Value = "convertible"
$ws.Cells.Item(5,3).Value = 0
$ws.Cells.Item(5,4).Formula = "'" + "2026-08-15"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = 110.5
$ws.Cells.Item(5,6).Value = "Convertible at `$2.50/share"
$ws.Cells.Item(5,7).Value = "0% Convertible Senior Notes due 2026, conversion price `$2.50"
$ws.Cells.Item(5,8).ClearContents()
$ws.Cells.Item(5,9).ClearContents()
$ws.Cells.Item(5,10).ClearContents()
$ws.Cells.Item(5,11).ClearContents()

# ---- Row 6 ----
$ws.Cells.Item(6,1).Value = "ABL Revolver Facility"
$ws.Cells.Item(6,2).Value = "secured"
$ws.Cells.Item(6,3).Value = 0
$ws.Cells.Item(6,4).Formula = "'" + "2026-12-31"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = 150
$ws.Cells.Item(6,6).Value = "Revolving credit"
$ws.Cells.Item(6,7).Value = "Asset-Based Lending Revolver, `$500M total capacity, `$150M drawn"
$ws.Cells.Item(6,8).ClearContents()
$ws.Cells.Item(6,9).ClearContents()
$ws.Cells.Item(6,10).ClearContents()
$ws.Cells.Item(6,11).ClearContents()

# ---- Row 7 ----
$ws.Cells.Item(7,1).Value = "Equipment Finance Leases"
$ws.Cells.Item(7,2).Value = "lease"
$ws.Cells.Item(7,3).Value = 8.2
$ws.Cells.Item(7,4).Formula = "'" + "2025-2030"
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).Value = 85.4
$ws.Cells.Item(7,6).ClearContents()
$ws.Cells.Item(7,7).Value = "Finance leases for theatre equipment and digital projectors"
$ws.Cells.Item(7,8).ClearContents()
$ws.Cells.Item(7,9).ClearContents()
$ws.Cells.Item(7,10).ClearContents()
$ws.Cells.Item(7,11).ClearContents()

# ---- Row 8 ----
$ws.Cells.Item(8,1).Value = "Odeon Term Loan"
$ws.Cells.Item(8,2).Value = "secured"
$ws.Cells.Item(8,3).Value = 0
$ws.Cells.Item(8,4).Formula = "'" + "2027-03-31"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Value = 245.8
$ws.Cells.Item(8,6).Value = "Callable with 2% premium"
$ws.Cells.Item(8,7).Value = "Odeon Cinemas UK term loan facility, secured by UK assets"
$ws.Cells.Item(8,8).ClearContents()
$ws.Cells.Item(8,9).ClearContents()
$ws.Cells.Item(8,10).ClearContents()
$ws.Cells.Item(8,11).ClearContents()

